$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the rounded timestamp already in row 3
$ws.Range("A3").Value = 45804.43432375

# New price entry (row 4), mirroring the layout/format of the existing rows
$ws.Range("A4").Value = 45804.43702679376
$ws.Range("A4").NumberFormat = $ws.Range("A3").NumberFormat
$ws.Range("B4").Value = "EVOWHEY PROTEIN"
$ws.Range("C4").Value = "2Kg"
$ws.Range("D4").Value = "37,90€"
